$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.173.09"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "2.508.41"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "494.38"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.66"
$ws.Range("E6").Value = "  +11.07%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +3.04%  "
$ws.Range("D9").Value = "2.530.43"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("E11").Value = "  +5.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  +4.13%  "
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "2.944.51"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").Value = "57.311.21"
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.36"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "2.524.27"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  +5.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.32"
$ws.Range("E20").Value = "  +5.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.45"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  +4.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.53"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "2.615.02"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.63"
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("D30").Value = "0.0₃0828"
$ws.Range("E30").Value = "  +6.51%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.68"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.53"
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.31"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  +3.24%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.84"
$ws.Range("E36").Value = "  +6.28%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.17"
$ws.Range("E37").Value = "  +5.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.887"
$ws.Range("E38").Value = "  +4.68%  "
$ws.Range("E39").Value = "  +9.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.40"
$ws.Range("E40").Value = "  +4.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  +4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.623"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0564"
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.94"
$ws.Range("E45").Value = "  +6.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "268.89"
$ws.Range("E46").Value = "  +5.00%  "
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.20"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.10"
$ws.Range("E50").Value = "  +6.19%  "
$ws.Range("D51").Value = "1.899.43"
$ws.Range("E51").Value = "  -1.86%  "
